# Fruta / hortaliza, semanal
#
# A new weekly observation is inserted at the top of the "Piña" data block
# (row 106), pushing the existing rows 106:193 down by one (to 107:194).
# Inserting the row shifts every following row's values/styles down
# automatically, so the record that falls off the bottom of the original
# block (old row 193) simply lands intact in the new row 194 - no explicit
# copy of the tail rows is needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 106 (rows 106:193 shift down to 107:194).
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with this week's record.
$ws.Cells.Item(106, 1).Value = 5
$ws.Cells.Item(106, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(106, 3).Value = "Maule"
$ws.Cells.Item(106, 4).Value = 44566
$ws.Cells.Item(106, 5).Value = 7
$ws.Cells.Item(106, 6).Value = "Fruta"
$ws.Cells.Item(106, 7).Value = 100108
$ws.Cells.Item(106, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(106, 9).Value = 100108005
$ws.Cells.Item(106, 10).Value = "Piña"
$ws.Cells.Item(106, 11).Value = "Caramelo"
$ws.Cells.Item(106, 12).Value = "Segunda"
$ws.Cells.Item(106, 13).Value = 270
$ws.Cells.Item(106, 14).Value = 14000
$ws.Cells.Item(106, 15).Value = 14000
$ws.Cells.Item(106, 16).Value = 14000
$ws.Cells.Item(106, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(106, 18).Value = "Ecuador"
$ws.Cells.Item(106, 19).Value = 1000
$ws.Cells.Item(106, 20).Value = 14
